$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 578 entirely (post removed), causing subsequent rows to shift up by one.
$ws.Rows.Item(578).Delete()
